# Update "想去人数" (attendee count) figures in the 江西-漫展信息 workbook.
# Sheet "展览" holds one set of rows; sheet "全部类型" mirrors the same
# events at different row numbers, so both need to be updated.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 2864
$wsExhibit.Range("F7").Value = 255
$wsExhibit.Range("F10").Value = 91
$wsExhibit.Range("F11").Value = 119
$wsExhibit.Range("F12").Value = 36
$wsExhibit.Range("F13").Value = 2698
$wsExhibit.Range("F14").Value = 931

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 2864
$wsAll.Range("F8").Value = 255
$wsAll.Range("F12").Value = 91
$wsAll.Range("F13").Value = 119
$wsAll.Range("F14").Value = 36
$wsAll.Range("F15").Value = 2698
$wsAll.Range("F16").Value = 931
